$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new mapped column entry "opm_id" in cell B9 (matches row 9,
# which previously only had the style placeholder in column C)
$ws.Range("B9").Value = "opm_id"
